# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (s_vals) values computed/regenerated for rows 2-14, replacing the
# previous "Strike#"-derived values in column G.
$sVals = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 3
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 0
}

foreach ($row in $sVals.Keys) {
    $ws.Range("G$row").Value = $sVals[$row]
}
